$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set "Who is doing" (column B) for question rows 6 and 7 (sheet rows 7 and 8) to "Tomek"
$ws.Range("B7").Value = "Tomek"
$ws.Range("B8").Value = "Tomek"

# Update "Status" (column C) for those same rows to "Done"
$ws.Range("C7").Value = "Done"
$ws.Range("C8").Value = "Done"

# Move the active selection to B8
$ws.Range("B8").Select()
